# Replace the "<<judgeRecital>>" / "<<writtenOrder>>" placeholder paragraphs
# with the new recital/order sentence, merging the two placeholder
# paragraphs (and the blank paragraph between them) into a single
# paragraph.

$d = $word.ActiveDocument

# Locate the two placeholder paragraphs by their field-code text.
$judgeIdx = -1
$writtenIdx = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($judgeIdx -eq -1 -and $t -like "*<<judgeRecital>>*") { $judgeIdx = $i }
    if ($writtenIdx -eq -1 -and $t -like "*<<writtenOrder>>*") { $writtenIdx = $i }
}

if ($judgeIdx -eq -1 -or $writtenIdx -eq -1) {
    throw "Could not locate judgeRecital/writtenOrder placeholder paragraphs"
}

# New OOXML for the merged paragraph's content (runs + proofErr spell
# markers around the merge-field names, matching the document's existing
# convention for <<mergeField>> placeholders).
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Upon the application </w:t></w:r><w:r><w:t>of &lt;&lt;</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>a</w:t></w:r><w:r><w:t>pplicant</w:t></w:r><w:r><w:t>N</w:t></w:r><w:r><w:t>ame</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>&gt;&gt; dated</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>&lt;&lt;</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>a</w:t></w:r><w:r><w:t>pplication</w:t></w:r><w:r><w:t>D</w:t></w:r><w:r><w:t>ate</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">&gt;&gt; </w:t></w:r><w:r><w:t>and upon considering the information provided by the parties</w:t></w:r><w:r><w:t>, the court has ordered written representations from both the applicant and the respondent.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# Replace the *content* of the judgeRecital paragraph (but not its own
# paragraph mark) so the paragraph keeps its original identity/formatting.
$judgePara = $d.Paragraphs.Item($judgeIdx)
$contentRange = $d.Range($judgePara.Range.Start, $judgePara.Range.End - 1)
$contentRange.InsertXML($xml)

# Remove the (now stale) writtenOrder paragraph entirely - its text has
# been folded into the judgeRecital paragraph above.
$writtenPara = $d.Paragraphs.Item($writtenIdx)
$writtenPara.Range.Delete()

# Remove the blank paragraph that used to sit between judgeRecital and
# writtenOrder, leaving just the single blank paragraph that originally
# followed writtenOrder.
$blankPara = $d.Paragraphs.Item($judgeIdx + 1)
$blankPara.Range.Delete()
